$wb = $excel.ActiveWorkbook

# Insert a new worksheet "MarkerTagIssues" immediately before "Notes"
$notes = $wb.Worksheets.Item("Notes")
$ws = $wb.Worksheets.Add($notes)
$ws.Name = "MarkerTagIssues"

# Write header row. The order in which *new* string literals are first
# written determines their position in the shared-strings table, so we
# write them in the same order the original authoring tool used
# (Site, TagNumber, Explanation, IssueEndDatetime, IssueStartDatetime)
# even though that isn't the left-to-right column order on the sheet.
$ws.Range("A1").Value = "Site"
$ws.Range("D1").Value = "TagNumber"
$ws.Range("E1").Value = "Explanation"
$ws.Range("C1").Value = "IssueEndDatetime"
$ws.Range("B1").Value = "IssueStartDatetime"
$ws.Range("D1").NumberFormat = "0"

# Row 2 - RB1 / tag 5394
$ws.Range("A2").Value = "RB1"
$ws.Range("B2").Value = 44201.541666666664
$ws.Range("C2").Value = 44229.536111111112
$ws.Range("D2").Value = 5394
$ws.Range("E2").Value = "If I remember correclty, it's possible that we didn't turn the antenna correctly/turn it back on after a tune. Resolved at next site visit"

# Row 3 - RB1 / tag 5394 (ongoing issue, end date is text, not a date)
$ws.Range("A3").Value = "RB1"
$ws.Range("B3").Value = 45059.11041666667
$ws.Range("E3").Value = "believe it ran out of batteries and we coulnd't get to it. Had some detections in october for some reason?"
$ws.Range("C3").Value = "Next site visit hopefully"
$ws.Range("D3").Value = 5394

# Row 4 - RB2 / tag 2102
$ws.Range("A4").Value = "RB2"
$ws.Range("B4").Value = 44139.497893518521
$ws.Range("C4").Value = 44168.577337962961
$ws.Range("D4").Value = 2102
$ws.Range("E4").Value = "not sure"

# Date/time formatting on the datetime cells (built-in format 22)
$ws.Range("B2:B4").NumberFormat = "m/d/yy h:mm"
$ws.Range("C2").NumberFormat = "m/d/yy h:mm"
$ws.Range("C4").NumberFormat = "m/d/yy h:mm"
# Integer formatting on the TagNumber column
$ws.Range("D2:D4").NumberFormat = "0"

# Column widths to fit the new content
$ws.Columns.Item(2).ColumnWidth = 13.833333333333334
$ws.Columns.Item(3).ColumnWidth = 19.333333333333332
$ws.Columns.Item(4).ColumnWidth = 12.833333333333334

# Page setup (portrait) to match the other data sheets
$ws.PageSetup.Orientation = 1

# Selection on the new sheet
$ws.Range("D6").Select() | Out-Null

# The new sheet becomes the active / visible tab (was AntennaMetadata before)
$ws.Activate()
